# Commit: "Add Inflation_contributions_graph_data.XLSX, Waterfall_graph_data.XLSX,
# and trimmed_graph_data.XLSX" -- appends a new data row (row 67) with one numeric
# value per column (A..AK) to the end of the existing table on Sheet 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 67 values, in column order A (1) .. AK (37).
$rowValues = @(
    0.00000430516198150763,  # A67
    0.0000242404473523366,  # B67
    0.000850384329134998,  # C67
    0.00000715645339718919,  # D67
    0.0000124005505310229,  # E67
    0,  # F67
    0.0000083191319076653,  # G67
    0.000000506793368285559,  # H67
    0.0000394616954346629,  # I67
    0.00000316081871442122,  # J67
    0.0000764253653139573,  # K67
    0.0000368583902290868,  # L67
    0,  # M67
    0.0000628511461420136,  # N67
    0.0000153462963170319,  # O67
    0.000120445029725793,  # P67
    0.000000829587308857147,  # Q67
    0.0000122660261059779,  # R67
    0.000169341473023423,  # S67
    -0.0000255058205752899,  # T67
    0.00000755419259070775,  # U67
    0.000000521273178808003,  # V67
    0.000000034318800873394,  # W67
    0.000197416947727098,  # X67
    0.0000932071336012287,  # Y67
    0.00173362265553749,  # Z67
    0.0000779773733267475,  # AA67
    0.00014055211861798,  # AB67
    0.000189752267900893,  # AC67
    0.00000739034476498405,  # AD67
    0,  # AE67
    0.00000114998050740886,  # AF67
    0.000125054825058275,  # AG67
    0.00000209689964667623,  # AH67
    0,  # AI67
    0.00000115568872143382,  # AJ67
    0.0000404340849016058  # AK67
)

$targetRow = 67
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item($targetRow, $col).Value = $rowValues[$col - 1]
}
